# Add a new row (row 5) to the "数组" (Array) sheet for the
# "Maximum Subarray" (最大子数组和) problem, as described by the commit
# "sub array with array".
#
# Column layout on this sheet: A=No. B=leetcode C=题目 D=解题方法
# E=解题关键词 F=时间复杂度 G=空间复杂度

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

$method = @"
1 dp[i]的长度等同于原数组长度，记录以原数组索引i以及之前的连续数组之和最大值
2 初始值是num[0]表示索引0以及之前的连续之和的最大值是nums[0]
3 迭代开始位置是1
4 如果dp[i-1]是负值，说明之前的先保存起来，要从nums[i]开始从新计算
5 如果dp[i-1]是正值，说明nums[i-1]有增益，就做累加
6 计算dp中的最大值
"@

$problem = @"
给定一个整数数组 nums ，找到一个具有最大和的连续子数组（子数组最少包含一个元素），返回其最大和。 
 示例: 
 输入: [-2,1,-3,4,-1,2,1,-5,4],
输出: 6
解释: 连续子数组 [4,-1,2,1] 的和最大，为 6。
 进阶: 
 如果你已经实现复杂度为 O(n) 的解法，尝试使用更为精妙的分治法求解。 
 Related Topics 数组 分治算法 动态规划
"@

$keywords = @"
动态规划
累加
最大值
"@

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 53

# Write D5 (method) before C5 (problem) so the new shared strings land at
# the same indices (161, 162, 163) as in the authored workbook.
$ws.Range("D5").Value = $method
$ws.Range("C5").Value = $problem
$ws.Range("E5").Value = $keywords

$ws.Range("F5").Value = "O(N)"
$ws.Range("G5").Value = "O(N)"

# Row 5 holds a lot of wrapped text; give it the same row height Excel
# would compute for similarly-sized wrapped content elsewhere in this
# workbook.
$ws.Rows.Item(5).RowHeight = 308

# Move the active selection to the newly-added method cell.
$ws.Range("D5").Select() | Out-Null
